$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 474864.16
$ws.Range("J17").Value = 603536.2
$ws.Range("L17").Value = 1810608.6
$ws.Range("N17").Value = -1810944.6

$ws.Range("H33").Value = 10260968
$ws.Range("I33").Value = 14922608
$ws.Range("K33").Value = 14922608
$ws.Range("M33").Value = -14922379

$ws.Range("H40").Value = 1175.2354
$ws.Range("I40").Value = 796.3333
$ws.Range("J40").Value = 1256.4286
$ws.Range("K40").Value = 796.3333
$ws.Range("L40").Value = 1256.4286
$ws.Range("M40").Value = -621.3333
$ws.Range("N40").Value = -1606.4286

$ws.Range("H51").Value = 83506664
$ws.Range("I51").Value = 207995.8
$ws.Range("K51").Value = 207995.8
$ws.Range("M51").Value = -207511.8

$ws.Range("H61").Value = 797.2222
$ws.Range("I61").Value = 459.375
$ws.Range("K61").Value = 1378.125
$ws.Range("M61").Value = -1206.125

$ws.Range("H62").Value = 3639.4
$ws.Range("I62").Value = 3639.4
$ws.Range("K62").Value = 3639.4
$ws.Range("M62").Value = -3015.4

$ws.Range("H65").Value = 3639.4
$ws.Range("I65").Value = 3639.4
$ws.Range("K65").Value = 18197
$ws.Range("M65").Value = -15077

$ws.Range("H88").Value = 1528
$ws.Range("I88").Value = 1450
$ws.Range("K88").Value = 1450
$ws.Range("M88").Value = -1044

$ws.Range("H91").Value = 1528
$ws.Range("I91").Value = 1450
$ws.Range("K91").Value = 1450
$ws.Range("M91").Value = -46

$ws.Range("H106").Value = 1694.7368
$ws.Range("I106").Value = 1455.5454
$ws.Range("J106").Value = 3273.4
$ws.Range("K106").Value = 1455.5454
$ws.Range("L106").Value = 3273.4
$ws.Range("M106").Value = -824.5454
$ws.Range("N106").Value = -4535.4

$ws.Range("H112").Value = 34658.844
$ws.Range("I112").Value = 2197.6
$ws.Range("K112").Value = 6592.799999999999
$ws.Range("M112").Value = -5484.799999999999

$ws.Range("H116").Value = 3920.2727
$ws.Range("I116").Value = 4078
$ws.Range("K116").Value = 4078
$ws.Range("M116").Value = -636

$ws.Range("H132").Value = 39420.035
$ws.Range("I132").Value = 43973.207
$ws.Range("K132").Value = 131919.621
$ws.Range("M132").Value = -129389.621

$ws.Range("H137").Value = 50000800
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 580.1
$ws.Range("I4").Value = 486.57144
$ws.Range("K4").Value = 486.57144
$ws.Range("M4").Value = -370.57144

$ws.Range("H32").Value = 168208.39
$ws.Range("I32").Value = 173405
$ws.Range("K32").Value = 173405
$ws.Range("M32").Value = -173118

$ws.Range("H45").Value = 57957.11
$ws.Range("I45").Value = 79195.16
$ws.Range("K45").Value = 79195.16
$ws.Range("M45").Value = -78818.16

$ws.Range("H74").Value = 407070.94
$ws.Range("I74").Value = 1238.8206
$ws.Range("J74").Value = 1396286.8
$ws.Range("K74").Value = 1238.8206
$ws.Range("L74").Value = 1396286.8
$ws.Range("M74").Value = -364.8206
$ws.Range("N74").Value = -1398034.8

$ws.Range("H77").Value = 407070.94
$ws.Range("I77").Value = 1238.8206
$ws.Range("J77").Value = 1396286.8
$ws.Range("K77").Value = 6194.103
$ws.Range("L77").Value = 6981434
$ws.Range("M77").Value = -1826.103
$ws.Range("N77").Value = -6990170

$ws.Range("H102").Value = 1780.4736
$ws.Range("I102").Value = 1709.8667
$ws.Range("K102").Value = 1709.8667
$ws.Range("M102").Value = -87.86670000000004

$ws.Range("H110").Value = 1530
$ws.Range("I110").Value = 1530
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1530
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 515
$ws.Range("N110").ClearContents()

$ws.Range("H119").Value = 58674.25
$ws.Range("J119").Value = 58674.25
$ws.Range("L119").Value = 58674.25
$ws.Range("N119").Value = -68350.25

$ws.Range("H122").Value = 1045.2
$ws.Range("I122").Value = 1045.2
$ws.Range("K122").Value = 3135.6
$ws.Range("M122").Value = -685.6000000000004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3656.2856
$ws.Range("I86").Value = 1539
$ws.Range("K86").Value = 1539
$ws.Range("M86").Value = -416

$ws.Range("H89").Value = 3656.2856
$ws.Range("I89").Value = 1539
$ws.Range("K89").Value = 7695
$ws.Range("M89").Value = -2079

$ws.Range("H99").Value = 13135.909
$ws.Range("I99").Value = 21999.334
$ws.Range("K99").Value = 21999.334
$ws.Range("M99").Value = -20501.334

$ws.Range("H105").Value = 13418.5
$ws.Range("I105").Value = 15741
$ws.Range("K105").Value = 15741
$ws.Range("M105").Value = -13994

$ws.Range("H107").Value = 9053.192999999999
$ws.Range("I107").Value = 10164.385
$ws.Range("J107").Value = 3275
$ws.Range("K107").Value = 10164.385
$ws.Range("L107").Value = 3275
$ws.Range("M107").Value = -8244.385
$ws.Range("N107").Value = -7115

$ws.Range("H134").Value = 22502176
$ws.Range("I134").Value = 1935.7059
$ws.Range("K134").Value = 5807.1177
$ws.Range("M134").Value = -3272.1177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1348.2858
$ws.Range("I22").Value = 1223
$ws.Range("K22").Value = 1223
$ws.Range("M22").Value = -873

$ws.Range("H58").Value = 2089.5862
$ws.Range("J58").Value = 1655.9286
$ws.Range("L58").Value = 1655.9286
$ws.Range("N58").Value = -2061.9286

$ws.Range("H62").Value = 10021.889
$ws.Range("I62").Value = 10028.286
$ws.Range("K62").Value = 10028.286
$ws.Range("M62").Value = -9404.286

$ws.Range("H65").Value = 10021.889
$ws.Range("I65").Value = 10028.286
$ws.Range("K65").Value = 50141.43
$ws.Range("M65").Value = -47021.43

$ws.Range("H134").Value = 1450.0625
$ws.Range("I134").Value = 1312.2106
$ws.Range("K134").Value = 3936.6318
$ws.Range("M134").Value = -1401.6318

$ws.Range("H136").Value = 2089.5862
$ws.Range("J136").Value = 1655.9286
$ws.Range("L136").Value = 4967.7858
$ws.Range("N136").Value = -10067.7858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 101.818184
$ws.Range("I2").Value = 148.28572
$ws.Range("J2").Value = 20.5
$ws.Range("K2").Value = 889.71432
$ws.Range("L2").Value = 123
$ws.Range("M2").Value = -776.71432
$ws.Range("N2").Value = -349

$ws.Range("H4").Value = 4720411.5
$ws.Range("J4").Value = 400639.8
$ws.Range("L4").Value = 1201919.4
$ws.Range("N4").Value = -1202143.4

$ws.Range("H35").Value = 467.5
$ws.Range("I35").Value = 495
$ws.Range("J35").Value = 440
$ws.Range("K35").Value = 1485
$ws.Range("L35").Value = 1320
$ws.Range("M35").Value = -1197
$ws.Range("N35").Value = -1896

$ws.Range("H56").Value = 9621071
$ws.Range("I56").Value = 9621071
$ws.Range("K56").Value = 9621071
$ws.Range("M56").Value = -9620541

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2759
$ws.Range("I102").Value = 2361.25
$ws.Range("K102").Value = 2361.25
$ws.Range("M102").Value = -739.25

$ws.Range("H132").Value = 922281.8
$ws.Range("I132").Value = 7806.1577
$ws.Range("K132").Value = 23418.4731
$ws.Range("M132").Value = -20888.4731

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 9800
$ws.Range("I29").Value = 9800
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 9800
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("M29").Value = -9505

$ws.Range("H122").Value = 3403.0908
$ws.Range("I122").Value = 2974.1072
$ws.Range("K122").Value = 8922.321599999999
$ws.Range("M122").Value = -6472.321599999999

$ws.Range("H132").Value = 3080.9707
$ws.Range("I132").Value = 2681.5454
$ws.Range("K132").Value = 8044.6362
$ws.Range("M132").Value = -5514.6362

$ws.Range("H136").Value = 2108.4768
$ws.Range("J136").Value = 1779.875
$ws.Range("L136").Value = 5339.625
$ws.Range("N136").Value = -10439.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 38000
$ws.Range("J119").Value = 38000
$ws.Range("L119").Value = 38000
$ws.Range("N119").Value = -47676

$ws.Range("H136").Value = 23171.436
$ws.Range("I136").Value = 29825.766
$ws.Range("K136").Value = 89477.298
$ws.Range("M136").Value = -86927.298
